# NIT-9007399230.xlsx - "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The detail table (rows 16-23) lists one row per worker/mora-period. The
# refreshed statement only keeps the "1910" period for the first two
# workers (WILMER -> now row 16, MARIA ZULEIDA -> now row 17) and drops the
# rows for the other periods (2507/2505) as well as the third worker
# (YURANIS), whose rows are removed entirely. The footer block (signature
# lines) simply shifts up to follow the shorter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (MARIA ZULEIDA / period 1910) becomes the new last row of the
# detail table once the rows below/around it are gone, so first clone the
# "closing" border formatting that currently lives on row 23 onto row 21.
$ws.Range("B23:J23").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# Drop the obsolete detail rows. Deleting from the bottom up keeps the
# row numbers of everything still above the cursor stable:
#   23  -> YURANIS SAYONARA LEON ATENCIO, period 2505
#   22  -> YURANIS SAYONARA LEON ATENCIO, period 2507
#   20  -> MARIA ZULEIDA RODRIGUEZ ALMEIDA, period 2505
#   19  -> MARIA ZULEIDA RODRIGUEZ ALMEIDA, period 2507  (period 1910 in row 21 is kept -> becomes row 17)
#   17  -> WILMER DE LA CANDELA LEON ATENCIO, period 2505
#   16  -> WILMER DE LA CANDELA LEON ATENCIO, period 2507 (period 1910 in row 18 is kept -> becomes row 16)
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(16).Delete()

# Refresh the summary figures for the updated statement.
$ws.Range("E11").Value = 66250
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1
